# "Generate Report for Handback" -- mark the two localized files as handed
# back (in sync with en-US), and record the "Latest Target File" (md) /
# "Latest Handback File" (xlf) hyperlinks + refreshed "Latest Handback
# DateTime" stamps on both the zh-cn and de-de status sheets.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# Cornflower blue (#6495ED), same ink used by the sheet's existing HyperLink
# cell style, expressed as the BGR-packed long the Excel object model wants.
$hyperlinkColor = 15570276

function Set-HandbackRow($ws, $row, $mdName, $mdUrl, $xlfName, $xlfUrl, $handbackStamp) {
    # B: Status -> handed back
    $ws.Cells.Item($row, 2).Value = $statusText

    # E: Latest Target File (the source .md, same target as column A's link)
    $eCell = $ws.Cells.Item($row, 5)
    $ws.Hyperlinks.Add($eCell, $mdUrl, "", "", $mdName) | Out-Null
    $eCell.Font.Name = "Calibri"
    $eCell.Font.Size = 11
    $eCell.Font.Underline = 2
    $eCell.Font.Color = $hyperlinkColor

    # F: Latest Handback File (the translated .xlf, same target as column C's link)
    $fCell = $ws.Cells.Item($row, 6)
    $ws.Hyperlinks.Add($fCell, $xlfUrl, "", "", $xlfName) | Out-Null
    $fCell.Font.Name = "Calibri"
    $fCell.Font.Size = 11
    $fCell.Font.Underline = 2
    $fCell.Font.Color = $hyperlinkColor

    # G: Latest Handback DateTime -> refreshed timestamp (plain text, like the
    # existing cells in this column)
    $ws.Cells.Item($row, 7).Value = $handbackStamp
}

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

Set-HandbackRow $wsZh 2 `
    "1d7ee5ab-25dc-427c-9ba7-d760e14e9d1b.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/c9241b1c07fa6c011c57177deb357db5419b8e84/e2e/1d7ee5ab-25dc-427c-9ba7-d760e14e9d1b.md" `
    "1d7ee5ab-25dc-427c-9ba7-d760e14e9d1b.63194497c13f7696af8acce6728b295184dea86a.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c80c218dd6f93e33d0b2eaaef4bf77ea962627b0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/1d7ee5ab-25dc-427c-9ba7-d760e14e9d1b.63194497c13f7696af8acce6728b295184dea86a.zh-cn.xlf" `
    "2016-03-07 05:11:38"

Set-HandbackRow $wsZh 3 `
    "29f983bc-ed00-4945-b663-5fc111d2269f.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/c9241b1c07fa6c011c57177deb357db5419b8e84/e2e/29f983bc-ed00-4945-b663-5fc111d2269f.md" `
    "29f983bc-ed00-4945-b663-5fc111d2269f.003e98a0eeedcc1d9f6f18cc32f67ddd27321de5.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c80c218dd6f93e33d0b2eaaef4bf77ea962627b0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/29f983bc-ed00-4945-b663-5fc111d2269f.003e98a0eeedcc1d9f6f18cc32f67ddd27321de5.zh-cn.xlf" `
    "2016-03-07 05:11:38"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

Set-HandbackRow $wsDe 2 `
    "1d7ee5ab-25dc-427c-9ba7-d760e14e9d1b.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/c9241b1c07fa6c011c57177deb357db5419b8e84/e2e/1d7ee5ab-25dc-427c-9ba7-d760e14e9d1b.md" `
    "1d7ee5ab-25dc-427c-9ba7-d760e14e9d1b.63194497c13f7696af8acce6728b295184dea86a.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b12d49e9f18fc84c0584d45869d6d5b217b56b46/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/1d7ee5ab-25dc-427c-9ba7-d760e14e9d1b.63194497c13f7696af8acce6728b295184dea86a.de-de.xlf" `
    "2016-03-07 05:11:57"

Set-HandbackRow $wsDe 3 `
    "29f983bc-ed00-4945-b663-5fc111d2269f.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/c9241b1c07fa6c011c57177deb357db5419b8e84/e2e/29f983bc-ed00-4945-b663-5fc111d2269f.md" `
    "29f983bc-ed00-4945-b663-5fc111d2269f.003e98a0eeedcc1d9f6f18cc32f67ddd27321de5.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b12d49e9f18fc84c0584d45869d6d5b217b56b46/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/29f983bc-ed00-4945-b663-5fc111d2269f.003e98a0eeedcc1d9f6f18cc32f67ddd27321de5.de-de.xlf" `
    "2016-03-07 05:11:57"
